$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

Set-TextValue $ws 'D2' '67.661.41'
$ws.Range('E2').Value = '  +0.77%  '
Set-TextValue $ws 'D3' '2.481.01'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws 'D5' '586.80'
$ws.Range('E5').Value = '  +0.23%  '
Set-TextValue $ws 'D6' '174.54'
$ws.Range('E6').Value = '  +1.33%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('E9').Value = '  +4.21%  '
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('E12').Value = '  -0.21%  '
Set-TextValue $ws 'D13' '2.933.58'
$ws.Range('E13').Value = '  -0.77%  '
Set-TextValue $ws 'D14' '25.20'
$ws.Range('E14').Value = '  -1.31%  '
Set-TextValue $ws 'D15' '67.827.76'
$ws.Range('E15').Value = '  +1.16%  '
$ws.Range('E16').Value = '  -0.65%  '
Set-TextValue $ws 'D17' '2.474.81'
$ws.Range('E17').Value = '  -0.59%  '
Set-TextValue $ws 'D18' '7.39'
$ws.Range('E18').Value = '  -3.32%  '
Set-TextValue $ws 'D19' '10.77'
$ws.Range('E19').Value = '  -2.00%  '
Set-TextValue $ws 'D20' '346.63'
$ws.Range('E20').Value = '  -1.23%  '
$ws.Range('E21').Value = '  +1.40%  '
Set-TextValue $ws 'D22' '1.00'
$ws.Range('E22').Value = '  -0.09%  '
Set-TextValue $ws 'D23' '70.70'
$ws.Range('E23').Value = '  +2.51%  '
Set-TextValue $ws 'D24' '4.17'
$ws.Range('E24').Value = '  -1.25%  '
$ws.Range('E25').Value = '  -6.78%  '
Set-TextValue $ws 'D26' '8.77'
$ws.Range('E26').Value = '  -4.58%  '
Set-TextValue $ws 'D27' '2.603.62'
$ws.Range('E27').Value = '  +0.01%  '
Set-TextValue $ws 'D28' '0.998'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('E29').Value = '  -2.30%  '
Set-TextValue $ws 'D30' '7.70'
$ws.Range('E30').Value = '  +0.11%  '
Set-TextValue $ws 'D31' '492.94'
$ws.Range('E31').Value = '  -3.14%  '
$ws.Range('E32').Value = '  -0.28%  '
Set-TextValue $ws 'D33' '1.76'
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('E34').Value = '  +0.00%  '
Set-TextValue $ws 'D35' '164.32'
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('E36').Value = '  +1.44%  '
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('E40').Value = '  -3.45%  '
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('E42').Value = '  -1.75%  '
Set-TextValue $ws 'D43' '4.75'
$ws.Range('E43').Value = '  -1.44%  '
$ws.Range('E44').Value = '  -0.84%  '
Set-TextValue $ws 'D45' '147.34'
$ws.Range('E45').Value = '  +2.72%  '
Set-TextValue $ws 'D46' '3.51'
$ws.Range('E46').Value = '  +0.85%  '
Set-TextValue $ws 'D47' '0.509'
$ws.Range('E47').Value = '  -1.30%  '
$ws.Range('E48').Value = '  -4.09%  '
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('E51').Value = '  -1.31%  '
